$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 45.54, 7),
    @(3, 40.64, 8),
    @(4, 48.050000000000004, 4),
    @(5, 38.83, 4),
    @(6, 40.18, 9),
    @(7, 55.45, 3),
    @(8, 18.38, 13),
    @(9, 18.54, 11),
    @(10, 30.060000000000002, 11),
    @(11, 18.38, 8),
    @(12, 20.260000000000002, 10),
    @(13, 29.75, 6),
    @(14, 15.57, 16),
    @(15, 18.21, 15),
    @(16, 30.060000000000002, 8),
    @(17, 14.68, 12),
    @(18, 20.09, 8),
    @(19, 32.04, 7),
    @(20, 54.58, 0),
    @(21, 30.93, 11),
    @(22, 36.590000000000003, 8),
    @(23, 35.67, 2),
    @(24, 26.63, 5),
    @(25, 36.130000000000003, 13),
    @(26, 51.230000000000004, 3),
    @(27, 37.730000000000004, 7),
    @(28, 35.340000000000003, 13),
    @(29, 71.14, 6),
    @(30, 44.21, 5),
    @(31, 53.09, 6),
    @(32, 71.14, 2)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 6).Value = $row[2]
}

# Rows 28-32 originally used the "General" number format; the edit brings
# them in line with the rest of column F ("#,##0" style already used by F2:F27).
$ws.Range("F27").Copy()
$ws.Range("F28:F32").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D2:D32").Select()
